$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# 1. Version: "1.0" -> "1.1"  (row 1, col 4)
$cell = $tbl.Cell(1, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "1.1"

# 2. Date Last Updated: "5/31/2015" -> "6/1/2015" split across two runs "6/" + "1/2015" (row 3, col 4)
$cell = $tbl.Cell(3, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = ""
$r.Collapse(1)
$r.InsertAfter("6/")
$r.Collapse(0)
$r.InsertAfter("1/2015")

# 3. UML Diagram file name: "Loan_Grant_UML.vsdx" -> "UML v2.0.vsdx" (row 4, col 2)
$cell = $tbl.Cell(4, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "UML v2.0.vsdx"

# 4. Brief Description: split sentence into multiple runs (row 8, col 2)
$cell = $tbl.Cell(8, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = ""
$r.Collapse(1)
$r.InsertAfter("As a ")
$r.Collapse(0)
$r.InsertAfter("Registered User")
$r.Collapse(0)
$r.InsertAfter(", I want to be able to enter my")
$r.Collapse(0)
$r.InsertAfter(" credentials")
$r.Collapse(0)
$r.InsertAfter(", so I can ")
$r.Collapse(0)
$r.InsertAfter("access additional features")
$r.Collapse(0)
$r.InsertAfter(".")

$d.Save()
